$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels that used a "IdXxx?" naming convention to "XxxId?" (int-style naming)
$ws.Range("G12").Value = "GymId?"
$ws.Range("H12").Value = "TrMachId?"

# Column M of the Subscription table used to reference ProgWorkId?, now it holds a Time value
$ws.Range("M23").Value = "Time"
$ws.Range("M24").Value = 12
$ws.Range("M25").Value = 5

# Update the remembered selection to match the saved workbook state
$ws.Range("M26").Select()
